# Update NATMI LR-pair (Jag2-Notch1) results sheet with newly recomputed TPM-based values
# (per "update scripts wuth new tpm" commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 19.35876433333334
$ws.Range("H2").Value = 58.07629300000001
$ws.Range("I2").Value = 0.9707667559429034
$ws.Range("J2").Value = 0.9707667559429034
$ws.Range("M2").Value = 48.42420966666666
$ws.Range("N2").Value = 145.272629
$ws.Range("O2").Value = 0.6311762527593259
$ws.Range("P2").Value = 0.6311762527593258
$ws.Range("Q2").Value = 937.4328629649219
$ws.Range("R2").Value = 8436.895766684298
$ws.Range("S2").Value = 0.6127249233193688
$ws.Range("T2").Value = 0.6127249233193687

# Row 3
$ws.Range("G3").Value = 19.35876433333334
$ws.Range("H3").Value = 58.07629300000001
$ws.Range("I3").Value = 0.9707667559429034
$ws.Range("J3").Value = 0.9707667559429034
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("O3").Value = 0.08928392431779728
$ws.Range("P3").Value = 0.08928392431779726
$ws.Range("Q3").Value = 132.6058837354436
$ws.Range("R3").Value = 1193.452953618992
$ws.Range("S3").Value = 0.08667386556783976
$ws.Range("T3").Value = 0.08667386556783975

# Row 4
$ws.Range("G4").Value = 19.35876433333334
$ws.Range("H4").Value = 58.07629300000001
$ws.Range("I4").Value = 0.9707667559429034
$ws.Range("J4").Value = 0.9707667559429034
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2795398229228769
$ws.Range("P4").Value = 0.2795398229228769
$ws.Range("Q4").Value = 415.1769262067313
$ws.Range("R4").Value = 3736.592335860582
$ws.Range("S4").Value = 0.2713679670556949
$ws.Range("T4").Value = 0.2713679670556949

# Row 5
$ws.Range("H5").Value = 0.9049070000000001
$ws.Range("I5").Value = 0.01512585579145048
$ws.Range("J5").Value = 0.01512585579145048
$ws.Range("M5").Value = 48.42420966666666
$ws.Range("N5").Value = 145.272629
$ws.Range("O5").Value = 0.6311762527593259
$ws.Range("P5").Value = 0.6311762527593258
$ws.Range("Q5").Value = 14.60646876561144
$ws.Range("R5").Value = 131.458218890503
$ws.Range("S5").Value = 0.00954708097822566
$ws.Range("T5").Value = 0.009547080978225658

# Row 6
$ws.Range("H6").Value = 0.9049070000000001
$ws.Range("I6").Value = 0.01512585579145048
$ws.Range("J6").Value = 0.01512585579145048
$ws.Range("M6").Value = 6.849914666666667
$ws.Range("O6").Value = 0.08928392431779728
$ws.Range("P6").Value = 0.08928392431779726
$ws.Range("S6").Value = 0.00135049576372578
$ws.Range("T6").Value = 0.00135049576372578

# Row 7
$ws.Range("H7").Value = 0.9049070000000001
$ws.Range("I7").Value = 0.01512585579145048
$ws.Range("J7").Value = 0.01512585579145048
$ws.Range("N7").Value = 64.33937399999999
$ws.Range("O7").Value = 0.2795398229228769
$ws.Range("P7").Value = 0.2795398229228769
$ws.Range("Q7").Value = 6.469016656468666
$ws.Range("R7").Value = 58.221149908218
$ws.Range("S7").Value = 0.004228279049499038
$ws.Range("T7").Value = 0.004228279049499038

# Row 8
$ws.Range("G8").Value = 0.2813256666666666
$ws.Range("H8").Value = 0.843977
$ws.Range("I8").Value = 0.01410738826564608
$ws.Range("J8").Value = 0.01410738826564608
$ws.Range("M8").Value = 48.42420966666666
$ws.Range("N8").Value = 145.272629
$ws.Range("O8").Value = 0.6311762527593259
$ws.Range("P8").Value = 0.6311762527593258
$ws.Range("Q8").Value = 13.62297306728144
$ws.Range("R8").Value = 122.606757605533
$ws.Range("S8").Value = 0.008904248461731379
$ws.Range("T8").Value = 0.008904248461731377

# Row 9
$ws.Range("G9").Value = 0.2813256666666666
$ws.Range("H9").Value = 0.843977
$ws.Range("I9").Value = 0.01410738826564608
$ws.Range("J9").Value = 0.01410738826564608
$ws.Range("M9").Value = 6.849914666666667
$ws.Range("O9").Value = 0.08928392431779728
$ws.Range("P9").Value = 0.08928392431779726
$ws.Range("S9").Value = 0.001259562986231726
$ws.Range("T9").Value = 0.001259562986231726

# Row 10
$ws.Range("G10").Value = 0.2813256666666666
$ws.Range("H10").Value = 0.843977
$ws.Range("I10").Value = 0.01410738826564608
$ws.Range("J10").Value = 0.01410738826564608
$ws.Range("N10").Value = 64.33937399999999
$ws.Range("O10").Value = 0.2795398229228769
$ws.Range("P10").Value = 0.2795398229228769
$ws.Range("Q10").Value = 6.033439094488665
$ws.Range("R10").Value = 54.30095185039799
$ws.Range("S10").Value = 0.003943576817682977
$ws.Range("T10").Value = 0.003943576817682977
